$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.701.72"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.847.91"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.75"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4275"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3622"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07315"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8708"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.73"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "1.852.65"
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.561"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.344"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06992"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.72"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008972"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.32"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "27.692.99"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.997"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.38"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("D24").Value = "2.072.83"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.987"
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.38"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.54"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "120.60"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.288"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.871"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08901"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7669"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.965"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.508"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("E35").Value = "  +3.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.004"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05435"
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.108"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01930"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.828"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1665"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5079"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.576"
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.424"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06551"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.48"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4656"
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.635"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.759"
$ws.Range("E51").Value = "  -0.17%  "
